$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D5").Value = "2016-02-26 07:10:20"
$wsZh.Range("G5").Value = "2016-02-26 07:11:33"

$wsDe.Range("D5").Value = "2016-02-26 07:10:38"
$wsDe.Range("G5").Value = "2016-02-26 07:12:00"
